$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 79243

# Row 4
$ws.Range("A4").Value = 130134267
$ws.Range("B4").Value = 99013
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = 'Knärot'
$ws.Range("G4").Value = 'Goodyera repens'
$ws.Range("H4").Value = '(L.) R. Br.'
$ws.Range("Q4").Value = 750666
$ws.Range("R4").Value = 7111136
$ws.Range("Y4").Formula = '="2025-08-29"'
$ws.Range("Y4").Copy()
$ws.Range("Y4").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("AA4").Formula = '="2025-08-29"'
$ws.Range("AA4").Copy()
$ws.Range("AA4").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("AX4").Value = 'Lisa Sandberg'

# Row 5
$ws.Range("A5").Value = 130134356
$ws.Range("B5").Value = 93095
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = 'Dropptaggsvamp'
$ws.Range("G5").Value = 'Hydnellum ferrugineum'
$ws.Range("H5").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q5").Value = 750732
$ws.Range("R5").Value = 7111314
$ws.Range("Y5").Formula = '="2025-10-07"'
$ws.Range("Y5").Copy()
$ws.Range("Y5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("AA5").Formula = '="2025-10-07"'
$ws.Range("AA5").Copy()
$ws.Range("AA5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("AX5").Value = 'Daniel Lussetti'

# Row 6
$ws.Range("B6").Value = 79243

# Row 7
$ws.Range("B7").Value = 93095

# Row 8
$ws.Range("B8").Value = 79243

# Row 9
$ws.Range("B9").Value = 79243

# Row 10
$ws.Range("B10").Value = 79000

# Row 11
$ws.Range("B11").Value = 79243

# Row 12
$ws.Range("B12").Value = 93107

# Row 13
$ws.Range("B13").Value = 79243

# Row 14
$ws.Range("B14").Value = 99013

# Row 15
$ws.Range("B15").Value = 79243

# Row 16
$ws.Range("B16").Value = 78646

# Row 17
$ws.Range("B17").Value = 79243

# Row 19
$ws.Range("B19").Value = 93133

# Row 20
$ws.Range("B20").Value = 93133
